$d = $word.ActiveDocument

function New-OpcPackageXml($bodyInnerXml) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
        $bodyInnerXml +
        '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# ---------------------------------------------------------------------------
# 1) Insert three new paragraphs (the "Accenture Labs" / speaker-name call
#    notes) right after the title paragraph, before the first blank
#    paragraph.
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$afterTitle = $titlePara.Range.End
$insertionPoint = $d.Range($afterTitle, $afterTitle)
$insertionPoint.InsertBefore("Accenture Labs`rIoannis Polykretis to Everyone (12:24 PM)`rUli Kremer`r")

# Re-fetch the two freshly created paragraphs that need spell-check
# (proofErr) run-splitting around the proper nouns, and rewrite their
# contents (minus the trailing paragraph mark) via InsertXML so the run /
# proofErr structure matches exactly.
$ioannisPara = $d.Paragraphs(3)
$ioannisRange = $d.Range($ioannisPara.Range.Start, $ioannisPara.Range.End - 1)
$ioannisXml = New-OpcPackageXml (
    '<w:p>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>Ioannis</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>Polykretis</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> to Everyone (12:24 PM)</w:t></w:r>' +
    '</w:p>'
)
$ioannisRange.InsertXML($ioannisXml)

$uliPara = $d.Paragraphs(4)
$uliRange = $d.Range($uliPara.Range.Start, $uliPara.Range.End - 1)
$uliXml = New-OpcPackageXml (
    '<w:p>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>Uli</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> Kremer</w:t></w:r>' +
    '</w:p>'
)
$uliRange.InsertXML($uliXml)

# ---------------------------------------------------------------------------
# 2) Add a lastRenderedPageBreak marker in front of the "gives you freedom
#    to explore ..." bullet line.
# ---------------------------------------------------------------------------
$findRange = $d.Content
$found = $findRange.Find.Execute(
    "gives you freedom to explore your own ideas whenever you feel competent to fly solo,",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $freedomRange = $d.Range($findRange.Start, $findRange.End)
    $freedomXml = New-OpcPackageXml (
        '<w:p><w:r><w:lastRenderedPageBreak/>' +
        '<w:t>gives you freedom to explore your own ideas whenever you feel competent to fly solo,</w:t>' +
        '</w:r></w:p>'
    )
    $freedomRange.InsertXML($freedomXml)
}

# ---------------------------------------------------------------------------
# 3) Move the lastRenderedPageBreak marker (and the run split point) earlier
#    in the "enthusiasm ... working together." paragraph: it used to sit
#    right before "inviting you to join the lab", it now sits right before
#    "admitted, feel free ...".
# ---------------------------------------------------------------------------
$oldFull = " enthusiasm. If you are invited for a visit, we will likely have a chance to talk about science. Tell me a project that you are most proud of or had most fun working on. Helping me to really understand one thing that you’ve done is likely to be more impressive than giving me a summary of many things. If you are admitted, feel free to drop by my office to discuss summer research and other rotation opportunities. We would always like to work with you for at least one rotation project before inviting you to join the lab, just so we get to know each other and only commit to the future if we remain mutually excited about working together."

$findRange2 = $d.Content
$found2 = $findRange2.Find.Execute($oldFull, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $enthusiasmRange = $d.Range($findRange2.Start, $findRange2.End)
    $newRun1 = " enthusiasm. If you are invited for a visit, we will likely have a chance to talk about science. Tell me a project that you are most proud of or had most fun working on. Helping me to really understand one thing that you’ve done is likely to be more impressive than giving me a summary of many things. If you are "
    $newRun2 = "admitted, feel free to drop by my office to discuss summer research and other rotation opportunities. We would always like to work with you for at least one rotation project before inviting you to join the lab, just so we get to know each other and only commit to the future if we remain mutually excited about working together."
    $enthusiasmXml = New-OpcPackageXml (
        '<w:p>' +
        '<w:r><w:t xml:space="preserve">' + $newRun1 + '</w:t></w:r>' +
        '<w:r><w:lastRenderedPageBreak/><w:t>' + $newRun2 + '</w:t></w:r>' +
        '</w:p>'
    )
    $enthusiasmRange.InsertXML($enthusiasmXml)
}
